$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "26.252.50"
$ws.Range("D3").Value = "1.606.10"
$ws.Range("E3").Value = "  +0.52%  "
Set-TextValue "D5" "212.71"
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("E6").Value = "  -0.02%  "
Set-TextValue "D7" "0.485"
$ws.Range("E7").Value = "  +0.04%  "
Set-TextValue "D8" "0.248"
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -0.25%  "
Set-TextValue "D10" "18.04"
$ws.Range("E10").Value = "  +0.74%  "
$ws.Range("E11").Value = "  -0.71%  "
$ws.Range("D12").Value = "1.832.45"
$ws.Range("E12").Value = "  +0.67%  "
$ws.Range("D13").Value = "1.602.84"
$ws.Range("E13").Value = "  -0.20%  "
Set-TextValue "D14" "4.02"
$ws.Range("E14").Value = "  +0.41%  "
Set-TextValue "D15" "0.513"
$ws.Range("E15").Value = "  +0.79%  "
$ws.Range("D16").Value = "26.256.46"
$ws.Range("E16").Value = "  +0.65%  "
Set-TextValue "D17" "62.07"
$ws.Range("E17").Value = "  +2.73%  "
$ws.Range("D18").Value = "0.0₃0727"
$ws.Range("E18").Value = "  +0.85%  "
$ws.Range("E19").Value = "  -0.11%  "
Set-TextValue "D20" "201.71"
$ws.Range("E20").Value = "  -1.42%  "
Set-TextValue "D21" "4.27"
$ws.Range("E21").Value = "  +0.85%  "
Set-TextValue "D22" "9.30"
$ws.Range("E22").Value = "  -0.08%  "
Set-TextValue "D23" "6.00"
$ws.Range("E23").Value = "  +0.51%  "
$ws.Range("E24").Value = "  +2.71%  "
Set-TextValue "D25" "144.85"
$ws.Range("E25").Value = "  +2.04%  "
Set-TextValue "D26" "1.01"
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("E27").Value = "  -3.26%  "
Set-TextValue "D28" "15.19"
$ws.Range("E28").Value = "  -0.04%  "
Set-TextValue "D29" "6.55"
$ws.Range("E29").Value = "  +1.79%  "
$ws.Range("E30").Value = "  +5.03%  "
Set-TextValue "D31" "1.17"
$ws.Range("E31").Value = "  +0.62%  "
Set-TextValue "D32" "3.20"
$ws.Range("E32").Value = "  +2.71%  "
$ws.Range("E33").Value = "  -2.10%  "
Set-TextValue "D35" "1.48"
$ws.Range("E35").Value = "  +0.78%  "
$ws.Range("D36").Value = "1.162.54"
$ws.Range("E36").Value = "  +5.07%  "
Set-TextValue "D37" "0.0165"
$ws.Range("E37").Value = "  +1.15%  "
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("E40").Value = "  +1.31%  "
Set-TextValue "D41" "0.497"
$ws.Range("E41").Value = "  +0.81%  "
Set-TextValue "D42" "0.782"
$ws.Range("E42").Value = "  +0.65%  "
Set-TextValue "D43" "5.28"
$ws.Range("E43").Value = "  +3.91%  "
$ws.Range("D44").Value = "1.744.98"
$ws.Range("E44").Value = "  +0.64%  "
Set-TextValue "D45" "91.92"
$ws.Range("E45").Value = "  -0.75%  "
Set-TextValue "D46" "1.52"
$ws.Range("E46").Value = "  +0.91%  "
Set-TextValue "D47" "54.12"
$ws.Range("E47").Value = "  +1.49%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D48" "0.0506"
$ws.Range("E48").Value = "  +0.15%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₇0972"
$ws.Range("E49").Value = "  -5.15%  "
Set-TextValue "D50" "0.407"
$ws.Range("E50").Value = "  -0.46%  "
$ws.Range("E51").Value = "  +0.06%  "
